# Fruta / hortaliza, semanal
# Insert 8 new weekly records (week of 2022-01-10, serial 44578) above the
# existing row 420 in the "Femacal de La Calera - Cereza" sheet, pushing the
# old rows 420:434 down to 428:442.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 blank rows above row 420; existing data (old rows 420-434) shifts
# down to 428-442.
$ws.Rows("420:427").Insert()

# Columns that stay constant for every record in this block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"
$unidad    = "`$/bandeja 10 kilos"
$origen    = "Provincia de Curic" + [char]0xF3
$kgUnidad  = 10

# New records for the week of serial date 44578.
$newRows = @(
    @{ Row=420; Variedad="Brooks";      Calidad="Especial"; Volumen=120; Min=7000; Max=7000; Prom=7000; Kg=700 },
    @{ Row=421; Variedad="Brooks";      Calidad="Primera";  Volumen=110; Min=5000; Max=5000; Prom=5000; Kg=500 },
    @{ Row=422; Variedad="Brooks";      Calidad="Segunda";  Volumen=100; Min=4000; Max=4000; Prom=4000; Kg=400 },
    @{ Row=423; Variedad="Lapins";      Calidad="Especial"; Volumen=90;  Min=6000; Max=6000; Prom=6000; Kg=600 },
    @{ Row=424; Variedad="Lapins";      Calidad="Primera";  Volumen=85;  Min=5000; Max=5000; Prom=5000; Kg=500 },
    @{ Row=425; Variedad="Lapins";      Calidad="Segunda";  Volumen=80;  Min=4000; Max=4000; Prom=4000; Kg=400 },
    @{ Row=426; Variedad="Sweet Heart"; Calidad="Primera";  Volumen=80;  Min=5000; Max=5000; Prom=5000; Kg=500 },
    @{ Row=427; Variedad="Sweet Heart"; Calidad="Segunda";  Volumen=90;  Min=4000; Max=4000; Prom=4000; Kg=400 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = 44578
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $rec.Variedad
    $ws.Cells.Item($r, 12).Value = $rec.Calidad
    $ws.Cells.Item($r, 13).Value = $rec.Volumen
    $ws.Cells.Item($r, 14).Value = $rec.Min
    $ws.Cells.Item($r, 15).Value = $rec.Max
    $ws.Cells.Item($r, 16).Value = $rec.Prom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $rec.Kg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
